# Apply "more fixes for error messages" update to Bloom Charts workbook.
$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Sheet1")
$wsPitch = $wb.Worksheets.Item("pitch breakdown")
$wsPitcher = $wb.Worksheets.Item("pitcher breakdown")

# ---------------------------------------------------------------------------
# Sheet1: correct the first logged pitch and add four newly logged pitches
# ---------------------------------------------------------------------------
$wsMain.Range("B2").Value = "09/23/2023"
$wsMain.Range("C2").Value = "1"
$wsMain.Range("D2").Value = "90"
$wsMain.Range("E2").Value = "FB"
$wsMain.Range("F2").Value = "Foul Ball"

$wsMain.Range("A3").Value = "Andrew Armstrong"
$wsMain.Range("B3").Value = "09/23/2023"
$wsMain.Range("C3").Value = "2"
$wsMain.Range("D3").Value = "91"
$wsMain.Range("E3").Value = "FB"
$wsMain.Range("F3").Value = "Strike looking"

$wsMain.Range("A4").Value = "Andrew Armstrong"
$wsMain.Range("B4").Value = "09/23/2023"
$wsMain.Range("C4").Value = "3"
$wsMain.Range("D4").Value = "88"
$wsMain.Range("E4").Value = "FB"
$wsMain.Range("F4").Value = "Hit"

$wsMain.Range("A5").Value = "Andrew Armstrong"
$wsMain.Range("B5").Value = "09/23/2023"
$wsMain.Range("C5").Value = "4"
$wsMain.Range("D5").Value = "88"
$wsMain.Range("E5").Value = "SL"
$wsMain.Range("F5").Value = "Strike looking"

$wsMain.Range("A6").Value = "Andrew Armstrong"
$wsMain.Range("B6").Value = "09/23/2023"
$wsMain.Range("C6").Value = "5"
$wsMain.Range("D6").Value = "85"
$wsMain.Range("E6").Value = "SL"
$wsMain.Range("F6").Value = "Strike swing & miss"

# ---------------------------------------------------------------------------
# pitch breakdown: same new rows, with additional Strike/Swing/Free Bases info
# ---------------------------------------------------------------------------
$wsPitch.Range("B2").Value = "09/23/2023"
$wsPitch.Range("F2").Value = "Foul Ball"
$wsPitch.Range("G2").Value = "Strike"
$wsPitch.Range("H2").Value = "Swing contact"
$wsPitch.Range("I2").Value = "nothing"

$wsPitch.Range("A3").Value = "Andrew Armstrong"
$wsPitch.Range("B3").Value = "09/23/2023"
$wsPitch.Range("C3").Value = "2"
$wsPitch.Range("D3").Value = "91"
$wsPitch.Range("E3").Value = "FB"
$wsPitch.Range("F3").Value = "Strike looking"
$wsPitch.Range("G3").Value = "Strike"
$wsPitch.Range("H3").Value = "No swing"
$wsPitch.Range("I3").Value = "nothing"

$wsPitch.Range("A4").Value = "Andrew Armstrong"
$wsPitch.Range("B4").Value = "09/23/2023"
$wsPitch.Range("C4").Value = "3"
$wsPitch.Range("D4").Value = "88"
$wsPitch.Range("E4").Value = "FB"
$wsPitch.Range("F4").Value = "Hit"
$wsPitch.Range("G4").Value = "Strike"
$wsPitch.Range("H4").Value = "Swing contact"
$wsPitch.Range("I4").Value = "not free base"

$wsPitch.Range("A5").Value = "Andrew Armstrong"
$wsPitch.Range("B5").Value = "09/23/2023"
$wsPitch.Range("C5").Value = "4"
$wsPitch.Range("D5").Value = "88"
$wsPitch.Range("E5").Value = "SL"
$wsPitch.Range("F5").Value = "Strike looking"
$wsPitch.Range("G5").Value = "Strike"
$wsPitch.Range("H5").Value = "No swing"
$wsPitch.Range("I5").Value = "nothing"

$wsPitch.Range("A6").Value = "Andrew Armstrong"
$wsPitch.Range("B6").Value = "09/23/2023"
$wsPitch.Range("C6").Value = "5"
$wsPitch.Range("D6").Value = "85"
$wsPitch.Range("E6").Value = "SL"
$wsPitch.Range("F6").Value = "Strike swing & miss"
$wsPitch.Range("G6").Value = "Strike"
$wsPitch.Range("H6").Value = "Swing no contact"
$wsPitch.Range("I6").Value = "nothing"

# ---------------------------------------------------------------------------
# pitcher breakdown: recalculated summary numbers for the pitcher
# ---------------------------------------------------------------------------
$wsPitcher.Range("B2").Value = 89.7
$wsPitcher.Range("C2").Value = 91
$wsPitcher.Range("D2").Value = 1
$wsPitcher.Range("E2").Value = 0.3
$wsPitcher.Range("F2").Value = 0.6
$wsPitcher.Range("G2").Value = 0.3
$wsPitcher.Range("H2").Value = 1
$wsPitcher.Range("I2").Value = 0
